$d = $word.ActiveDocument
$t = $d.Tables(1)
$t.Cell(1,1).Range.Find.Execute("34+43=", $true, $false, $false, $false, $false, $true, 1, $false, "92-27=", 1) | Out-Null
$t.Cell(1,2).Range.Find.Execute("20+52=", $true, $false, $false, $false, $false, $true, 1, $false, "14-3=", 1) | Out-Null
$t.Cell(1,3).Range.Find.Execute("17+68=", $true, $false, $false, $false, $false, $true, 1, $false, "73-16=", 1) | Out-Null
$t.Cell(1,4).Range.Find.Execute("89-43=", $true, $false, $false, $false, $false, $true, 1, $false, "80+14=", 1) | Out-Null
$t.Cell(1,5).Range.Find.Execute("64-55=", $true, $false, $false, $false, $false, $true, 1, $false, "40-22=", 1) | Out-Null
$t.Cell(2,1).Range.Find.Execute("42-1=", $true, $false, $false, $false, $false, $true, 1, $false, "30+22=", 1) | Out-Null
$t.Cell(2,2).Range.Find.Execute("21+17=", $true, $false, $false, $false, $false, $true, 1, $false, "15+25=", 1) | Out-Null
$t.Cell(2,3).Range.Find.Execute("47+27=", $true, $false, $false, $false, $false, $true, 1, $false, "79+13=", 1) | Out-Null
$t.Cell(2,4).Range.Find.Execute("16+11=", $true, $false, $false, $false, $false, $true, 1, $false, "7-3=", 1) | Out-Null
$t.Cell(2,5).Range.Find.Execute("53-39=", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=", 1) | Out-Null
$t.Cell(3,1).Range.Find.Execute("34-29=", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=", 1) | Out-Null
$t.Cell(3,2).Range.Find.Execute("43+11=", $true, $false, $false, $false, $false, $true, 1, $false, "86-13=", 1) | Out-Null
$t.Cell(3,3).Range.Find.Execute("7+61=", $true, $false, $false, $false, $false, $true, 1, $false, "20+47=", 1) | Out-Null
$t.Cell(3,4).Range.Find.Execute("71-31=", $true, $false, $false, $false, $false, $true, 1, $false, "1+72=", 1) | Out-Null
$t.Cell(3,5).Range.Find.Execute("80-19=", $true, $false, $false, $false, $false, $true, 1, $false, "66-35=", 1) | Out-Null
$t.Cell(4,1).Range.Find.Execute("10+1=", $true, $false, $false, $false, $false, $true, 1, $false, "23+69=", 1) | Out-Null
$t.Cell(4,2).Range.Find.Execute("9+45=", $true, $false, $false, $false, $false, $true, 1, $false, "67-4=", 1) | Out-Null
$t.Cell(4,3).Range.Find.Execute("89-53=", $true, $false, $false, $false, $false, $true, 1, $false, "36+31=", 1) | Out-Null
$t.Cell(4,4).Range.Find.Execute("96-69=", $true, $false, $false, $false, $false, $true, 1, $false, "54-14=", 1) | Out-Null
$t.Cell(4,5).Range.Find.Execute("73-61=", $true, $false, $false, $false, $false, $true, 1, $false, "6+25=", 1) | Out-Null
$t.Cell(5,1).Range.Find.Execute("17-3=", $true, $false, $false, $false, $false, $true, 1, $false, "1+71=", 1) | Out-Null
$t.Cell(5,2).Range.Find.Execute("74-72=", $true, $false, $false, $false, $false, $true, 1, $false, "9+61=", 1) | Out-Null
$t.Cell(5,3).Range.Find.Execute("10+66=", $true, $false, $false, $false, $false, $true, 1, $false, "90-11=", 1) | Out-Null
$t.Cell(5,4).Range.Find.Execute("36+58=", $true, $false, $false, $false, $false, $true, 1, $false, "72-53=", 1) | Out-Null
$t.Cell(5,5).Range.Find.Execute("79-53=", $true, $false, $false, $false, $false, $true, 1, $false, "80-3=", 1) | Out-Null
$t.Cell(6,1).Range.Find.Execute("20+41=", $true, $false, $false, $false, $false, $true, 1, $false, "46+48=", 1) | Out-Null
$t.Cell(6,2).Range.Find.Execute("46+11=", $true, $false, $false, $false, $false, $true, 1, $false, "87+10=", 1) | Out-Null
$t.Cell(6,3).Range.Find.Execute("89-61=", $true, $false, $false, $false, $false, $true, 1, $false, "12-6=", 1) | Out-Null
$t.Cell(6,4).Range.Find.Execute("99-91=", $true, $false, $false, $false, $false, $true, 1, $false, "19+13=", 1) | Out-Null
$t.Cell(6,5).Range.Find.Execute("84-76=", $true, $false, $false, $false, $false, $true, 1, $false, "90-76=", 1) | Out-Null
$t.Cell(7,1).Range.Find.Execute("10+12=", $true, $false, $false, $false, $false, $true, 1, $false, "22-14=", 1) | Out-Null
$t.Cell(7,2).Range.Find.Execute("2+52=", $true, $false, $false, $false, $false, $true, 1, $false, "40+35=", 1) | Out-Null
$t.Cell(7,3).Range.Find.Execute("48-38=", $true, $false, $false, $false, $false, $true, 1, $false, "68-46=", 1) | Out-Null
$t.Cell(7,4).Range.Find.Execute("85-46=", $true, $false, $false, $false, $false, $true, 1, $false, "24+72=", 1) | Out-Null
$t.Cell(7,5).Range.Find.Execute("47-14=", $true, $false, $false, $false, $false, $true, 1, $false, "14+46=", 1) | Out-Null
$t.Cell(8,1).Range.Find.Execute("16+26=", $true, $false, $false, $false, $false, $true, 1, $false, "93-21=", 1) | Out-Null
$t.Cell(8,2).Range.Find.Execute("32-22=", $true, $false, $false, $false, $false, $true, 1, $false, "63+25=", 1) | Out-Null
$t.Cell(8,3).Range.Find.Execute("37-7=", $true, $false, $false, $false, $false, $true, 1, $false, "84-2=", 1) | Out-Null
$t.Cell(8,4).Range.Find.Execute("49+10=", $true, $false, $false, $false, $false, $true, 1, $false, "42+38=", 1) | Out-Null
$t.Cell(8,5).Range.Find.Execute("72-62=", $true, $false, $false, $false, $false, $true, 1, $false, "70+12=", 1) | Out-Null
$t.Cell(9,1).Range.Find.Execute("20-1=", $true, $false, $false, $false, $false, $true, 1, $false, "28+64=", 1) | Out-Null
$t.Cell(9,2).Range.Find.Execute("77-72=", $true, $false, $false, $false, $false, $true, 1, $false, "8+83=", 1) | Out-Null
$t.Cell(9,3).Range.Find.Execute("70-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67-0=", 1) | Out-Null
$t.Cell(9,4).Range.Find.Execute("68-29=", $true, $false, $false, $false, $false, $true, 1, $false, "50-43=", 1) | Out-Null
$t.Cell(9,5).Range.Find.Execute("17+17=", $true, $false, $false, $false, $false, $true, 1, $false, "57-54=", 1) | Out-Null
$t.Cell(10,1).Range.Find.Execute("21+2=", $true, $false, $false, $false, $false, $true, 1, $false, "1+95=", 1) | Out-Null
$t.Cell(10,2).Range.Find.Execute("53-50=", $true, $false, $false, $false, $false, $true, 1, $false, "43-0=", 1) | Out-Null
$t.Cell(10,3).Range.Find.Execute("89-31=", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=", 1) | Out-Null
$t.Cell(10,4).Range.Find.Execute("58-34=", $true, $false, $false, $false, $false, $true, 1, $false, "1+44=", 1) | Out-Null
$t.Cell(10,5).Range.Find.Execute("45-2=", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=", 1) | Out-Null
$t.Cell(11,1).Range.Find.Execute("50+38=", $true, $false, $false, $false, $false, $true, 1, $false, "2+48=", 1) | Out-Null
$t.Cell(11,2).Range.Find.Execute("88-25=", $true, $false, $false, $false, $false, $true, 1, $false, "85+5=", 1) | Out-Null
$t.Cell(11,3).Range.Find.Execute("49-39=", $true, $false, $false, $false, $false, $true, 1, $false, "91-81=", 1) | Out-Null
$t.Cell(11,4).Range.Find.Execute("53-28=", $true, $false, $false, $false, $false, $true, 1, $false, "67+10=", 1) | Out-Null
$t.Cell(11,5).Range.Find.Execute("82-30=", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=", 1) | Out-Null
$t.Cell(12,1).Range.Find.Execute("54-5=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 1) | Out-Null
$t.Cell(12,2).Range.Find.Execute("71+0=", $true, $false, $false, $false, $false, $true, 1, $false, "97-39=", 1) | Out-Null
$t.Cell(12,3).Range.Find.Execute("56-24=", $true, $false, $false, $false, $false, $true, 1, $false, "11-1=", 1) | Out-Null
$t.Cell(12,4).Range.Find.Execute("73-62=", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=", 1) | Out-Null
$t.Cell(12,5).Range.Find.Execute("34+0=", $true, $false, $false, $false, $false, $true, 1, $false, "3+44=", 1) | Out-Null
$t.Cell(13,1).Range.Find.Execute("1+32=", $true, $false, $false, $false, $false, $true, 1, $false, "43+8=", 1) | Out-Null
$t.Cell(13,2).Range.Find.Execute("80+13=", $true, $false, $false, $false, $false, $true, 1, $false, "7-5=", 1) | Out-Null
$t.Cell(13,3).Range.Find.Execute("27+35=", $true, $false, $false, $false, $false, $true, 1, $false, "13+20=", 1) | Out-Null
$t.Cell(13,4).Range.Find.Execute("57+9=", $true, $false, $false, $false, $false, $true, 1, $false, "43+53=", 1) | Out-Null
$t.Cell(13,5).Range.Find.Execute("14+44=", $true, $false, $false, $false, $false, $true, 1, $false, "10+27=", 1) | Out-Null
$t.Cell(14,1).Range.Find.Execute("32+40=", $true, $false, $false, $false, $false, $true, 1, $false, "71+26=", 1) | Out-Null
$t.Cell(14,2).Range.Find.Execute("7+53=", $true, $false, $false, $false, $false, $true, 1, $false, "65-33=", 1) | Out-Null
$t.Cell(14,3).Range.Find.Execute("74-33=", $true, $false, $false, $false, $false, $true, 1, $false, "12+76=", 1) | Out-Null
$t.Cell(14,4).Range.Find.Execute("60-33=", $true, $false, $false, $false, $false, $true, 1, $false, "19+53=", 1) | Out-Null
$t.Cell(14,5).Range.Find.Execute("43+38=", $true, $false, $false, $false, $false, $true, 1, $false, "47+50=", 1) | Out-Null
$t.Cell(15,1).Range.Find.Execute("44+41=", $true, $false, $false, $false, $false, $true, 1, $false, "91-14=", 1) | Out-Null
$t.Cell(15,2).Range.Find.Execute("28+36=", $true, $false, $false, $false, $false, $true, 1, $false, "18+78=", 1) | Out-Null
$t.Cell(15,3).Range.Find.Execute("8+87=", $true, $false, $false, $false, $false, $true, 1, $false, "53-24=", 1) | Out-Null
$t.Cell(15,4).Range.Find.Execute("39+23=", $true, $false, $false, $false, $false, $true, 1, $false, "32-20=", 1) | Out-Null
$t.Cell(15,5).Range.Find.Execute("14+68=", $true, $false, $false, $false, $false, $true, 1, $false, "27+65=", 1) | Out-Null
$t.Cell(16,1).Range.Find.Execute("96-94=", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=", 1) | Out-Null
$t.Cell(16,2).Range.Find.Execute("38+48=", $true, $false, $false, $false, $false, $true, 1, $false, "52+19=", 1) | Out-Null
$t.Cell(16,3).Range.Find.Execute("70+24=", $true, $false, $false, $false, $false, $true, 1, $false, "65-2=", 1) | Out-Null
$t.Cell(16,4).Range.Find.Execute("25+30=", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=", 1) | Out-Null
$t.Cell(16,5).Range.Find.Execute("29+48=", $true, $false, $false, $false, $false, $true, 1, $false, "79-41=", 1) | Out-Null
$t.Cell(17,1).Range.Find.Execute("42-26=", $true, $false, $false, $false, $false, $true, 1, $false, "63-10=", 1) | Out-Null
$t.Cell(17,2).Range.Find.Execute("22-8=", $true, $false, $false, $false, $false, $true, 1, $false, "25-5=", 1) | Out-Null
$t.Cell(17,3).Range.Find.Execute("46+30=", $true, $false, $false, $false, $false, $true, 1, $false, "13-1=", 1) | Out-Null
$t.Cell(17,4).Range.Find.Execute("1+42=", $true, $false, $false, $false, $false, $true, 1, $false, "40-13=", 1) | Out-Null
$t.Cell(17,5).Range.Find.Execute("67-20=", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=", 1) | Out-Null
$t.Cell(18,1).Range.Find.Execute("50+28=", $true, $false, $false, $false, $false, $true, 1, $false, "69+20=", 1) | Out-Null
$t.Cell(18,2).Range.Find.Execute("41-12=", $true, $false, $false, $false, $false, $true, 1, $false, "96-42=", 1) | Out-Null
$t.Cell(18,3).Range.Find.Execute("32+33=", $true, $false, $false, $false, $false, $true, 1, $false, "81-39=", 1) | Out-Null
$t.Cell(18,4).Range.Find.Execute("25-2=", $true, $false, $false, $false, $false, $true, 1, $false, "25+26=", 1) | Out-Null
$t.Cell(18,5).Range.Find.Execute("0+60=", $true, $false, $false, $false, $false, $true, 1, $false, "4+1=", 1) | Out-Null
$t.Cell(19,1).Range.Find.Execute("38+40=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 1) | Out-Null
$t.Cell(19,2).Range.Find.Execute("82-56=", $true, $false, $false, $false, $false, $true, 1, $false, "62-8=", 1) | Out-Null
$t.Cell(19,3).Range.Find.Execute("46+11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+37=", 1) | Out-Null
$t.Cell(19,4).Range.Find.Execute("91-87=", $true, $false, $false, $false, $false, $true, 1, $false, "39-9=", 1) | Out-Null
$t.Cell(19,5).Range.Find.Execute("92-24=", $true, $false, $false, $false, $false, $true, 1, $false, "84-49=", 1) | Out-Null
$t.Cell(20,1).Range.Find.Execute("10+86=", $true, $false, $false, $false, $false, $true, 1, $false, "48+13=", 1) | Out-Null
$t.Cell(20,2).Range.Find.Execute("82-15=", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=", 1) | Out-Null
$t.Cell(20,3).Range.Find.Execute("73-70=", $true, $false, $false, $false, $false, $true, 1, $false, "55-3=", 1) | Out-Null
$t.Cell(20,4).Range.Find.Execute("91-67=", $true, $false, $false, $false, $false, $true, 1, $false, "27+42=", 1) | Out-Null
$t.Cell(20,5).Range.Find.Execute("60+4=", $true, $false, $false, $false, $false, $true, 1, $false, "23+9=", 1) | Out-Null
